# Apply "#5: cash & deposit done" edit:
#  - sheet4 (現金/Cash): add columns E:K (property_category, category, date,
#    legislator_name, legislator_id, source_file, index) and turn row 1 into
#    a proper field-name header row.
#  - sheet5 (存款/Deposit): add columns G:M (property_category, category,
#    date, legislator_name, legislator_id, source_file, index) and turn row
#    1 into a proper field-name header row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "現金" (Cash) -> sheet4.xml
# ---------------------------------------------------------------------
$cash = $wb.Worksheets.Item("現金")

# Extend formatting of the header row (row 1, bold/border style) and the two
# data rows (rows 2-3, plain style) from the existing B:D columns into the
# newly-used E:K columns.
$cash.Range("B1:D1").Copy()
$cash.Range("E1:K1").PasteSpecial(-4122)
$cash.Range("B2:D2").Copy()
$cash.Range("E2:K2").PasteSpecial(-4122)
$cash.Range("B3:D3").Copy()
$cash.Range("E3:K3").PasteSpecial(-4122)

# Row 1 becomes a full field-name header (mirrors the schema used by the
# other "normal" output sheets such as 股票/基金受益憑證/保險).
$cash.Cells.Item(1,2).Value = "currency"
$cash.Cells.Item(1,3).Value = "owner"
$cash.Cells.Item(1,4).Value = "total"
$cash.Cells.Item(1,5).Value = "property_category"
$cash.Cells.Item(1,6).Value = "category"
$cash.Cells.Item(1,7).Value = "date"
$cash.Cells.Item(1,8).Value = "legislator_name"
$cash.Cells.Item(1,9).Value = "legislator_id"
$cash.Cells.Item(1,10).Value = "source_file"
$cash.Cells.Item(1,11).Value = "index"

# Row 2 (index 54)
$cash.Cells.Item(2,5).Value = "cash"
$cash.Cells.Item(2,6).Value = "normal"
$cash.Cells.Item(2,7).Value = "2013-11-22"
$cash.Cells.Item(2,8).Value = "吳育昇"
$cash.Cells.Item(2,9).Value = 1322
$cash.Cells.Item(2,10).Value = "tmp88481"
$cash.Cells.Item(2,11).Value = 54

# Row 3 (index 55)
$cash.Cells.Item(3,5).Value = "cash"
$cash.Cells.Item(3,6).Value = "normal"
$cash.Cells.Item(3,7).Value = "2013-11-22"
$cash.Cells.Item(3,8).Value = "吳育昇"
$cash.Cells.Item(3,9).Value = 1322
$cash.Cells.Item(3,10).Value = "tmp88481"
$cash.Cells.Item(3,11).Value = 55

# ---------------------------------------------------------------------
# Sheet "存款" (Deposit) -> sheet5.xml
# ---------------------------------------------------------------------
$dep = $wb.Worksheets.Item("存款")

# Extend formatting from B:F into the new G:M columns, for the header row
# and every one of the 20 data rows (rows 2-21).
$dep.Range("B1:F1").Copy()
$dep.Range("G1:M1").PasteSpecial(-4122)
for ($r = 2; $r -le 21; $r++) {
    $dep.Range("B$($r):F$($r)").Copy()
    $dep.Range("G$($r):M$($r)").PasteSpecial(-4122)
}

# Row 1 becomes a full field-name header.
$dep.Cells.Item(1,2).Value = "bank"
$dep.Cells.Item(1,3).Value = "deposit_type"
$dep.Cells.Item(1,4).Value = "currency"
$dep.Cells.Item(1,5).Value = "owner"
$dep.Cells.Item(1,6).Value = "total"
$dep.Cells.Item(1,7).Value = "property_category"
$dep.Cells.Item(1,8).Value = "category"
$dep.Cells.Item(1,9).Value = "date"
$dep.Cells.Item(1,10).Value = "legislator_name"
$dep.Cells.Item(1,11).Value = "legislator_id"
$dep.Cells.Item(1,12).Value = "source_file"
$dep.Cells.Item(1,13).Value = "index"

# Data rows 2-21: indices taken from column A (already present).
$indices = @(60,61,62,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80)
for ($r = 2; $r -le 21; $r++) {
    $idx = $indices[$r - 2]
    $dep.Cells.Item($r,7).Value = "deposit"
    $dep.Cells.Item($r,8).Value = "normal"
    $dep.Cells.Item($r,9).Value = "2013-11-22"
    $dep.Cells.Item($r,10).Value = "吳育昇"
    $dep.Cells.Item($r,11).Value = 1322
    $dep.Cells.Item($r,12).Value = "tmp88481"
    $dep.Cells.Item($r,13).Value = $idx
}
